$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Set column N (DOI) first for all new rows, to match shared-string insertion order
$ws.Range("N177").Value = '10.1063/1.5051514'
$ws.Range("N178").Value = '10.1063/1.5051514'
$ws.Range("N179").Value = '10.1063/1.5051514'
$ws.Range("N180").Value = '10.1063/1.5051514'
$ws.Range("N181").Value = '10.1063/1.5051514'
$ws.Range("N182").Value = '10.1063/1.5051514'
$ws.Range("N183").Value = '10.1063/1.5051514'
$ws.Range("N184").Value = '10.1063/1.5051514'
$ws.Range("N185").Value = '10.1063/1.5051514'

# Step 2: Set column B (Composition)
$ws.Range("B177").Value = 'ErGdHoLaTbY'
$ws.Range("B178").Value = 'DyGdHoLaTbY'
$ws.Range("B179").Value = 'DyErGdHoLuScTbY'
$ws.Range("B180").Value = 'ErGdHoLaTbY'
$ws.Range("B181").Value = 'DyGdHoLaTbY'
$ws.Range("B182").Value = 'DyErGdHoLuScTbY'
$ws.Range("B183").Value = 'ErGdHoLaTbY'
$ws.Range("B184").Value = 'DyGdHoLaTbY'
$ws.Range("B185").Value = 'DyErGdHoLuScTbY'

# Step 3: Set column E (Material Comment)
$ws.Range("E177").Value = 'strong O contamination concentrated in the grain boundaries'
$ws.Range("E178").Value = 'strong O contamination concentrated in the FCC (around 50% O)'
$ws.Range("E179").Value = 'strong O contamination concentrated in the FCC (around 50% O)'
$ws.Range("E180").Value = 'strong O contamination concentrated in the grain boundaries'
$ws.Range("E181").Value = 'strong O contamination concentrated in the FCC (around 50% O)'
$ws.Range("E182").Value = 'strong O contamination concentrated in the FCC (around 50% O)'
$ws.Range("E183").Value = 'strong O contamination concentrated in the grain boundaries'
$ws.Range("E184").Value = 'strong O contamination concentrated in the FCC (around 50% O)'
$ws.Range("E185").Value = 'strong O contamination concentrated in the FCC (around 50% O)'

# Step 4: Remaining text columns (reuse existing shared strings; order not significant)
$ws.Range("C177").Value = 'HCP'
$ws.Range("D177").Value = 'AAM'
$ws.Range("F177").Value = 'compressive yield stress'
$ws.Range("G177").Value = 'EXP'
$ws.Range("L177").Value = 'Pa'
$ws.Range("M177").Value = 'T4'
$ws.Range("C178").Value = 'HCP+FCC'
$ws.Range("D178").Value = 'AAM'
$ws.Range("F178").Value = 'compressive yield stress'
$ws.Range("G178").Value = 'EXP'
$ws.Range("L178").Value = 'Pa'
$ws.Range("M178").Value = 'T4'
$ws.Range("C179").Value = 'HCP+FCC'
$ws.Range("D179").Value = 'AAM'
$ws.Range("F179").Value = 'compressive yield stress'
$ws.Range("G179").Value = 'EXP'
$ws.Range("L179").Value = 'Pa'
$ws.Range("M179").Value = 'T4'
$ws.Range("C180").Value = 'HCP'
$ws.Range("D180").Value = 'AAM'
$ws.Range("F180").Value = 'compressive fracture strength'
$ws.Range("G180").Value = 'EXP'
$ws.Range("L180").Value = 'Pa'
$ws.Range("M180").Value = 'T4'
$ws.Range("C181").Value = 'HCP+FCC'
$ws.Range("D181").Value = 'AAM'
$ws.Range("F181").Value = 'compressive fracture strength'
$ws.Range("G181").Value = 'EXP'
$ws.Range("L181").Value = 'Pa'
$ws.Range("M181").Value = 'T4'
$ws.Range("C182").Value = 'HCP+FCC'
$ws.Range("D182").Value = 'AAM'
$ws.Range("F182").Value = 'compressive fracture strength'
$ws.Range("G182").Value = 'EXP'
$ws.Range("L182").Value = 'Pa'
$ws.Range("M182").Value = 'T4'
$ws.Range("C183").Value = 'HCP'
$ws.Range("D183").Value = 'AAM'
$ws.Range("F183").Value = 'compressive ductility'
$ws.Range("G183").Value = 'EXP'
$ws.Range("L183").Value = '%'
$ws.Range("M183").Value = 'T4'
$ws.Range("C184").Value = 'HCP+FCC'
$ws.Range("D184").Value = 'AAM'
$ws.Range("F184").Value = 'compressive ductility'
$ws.Range("G184").Value = 'EXP'
$ws.Range("L184").Value = '%'
$ws.Range("M184").Value = 'T4'
$ws.Range("C185").Value = 'HCP+FCC'
$ws.Range("D185").Value = 'AAM'
$ws.Range("F185").Value = 'compressive ductility'
$ws.Range("G185").Value = 'EXP'
$ws.Range("L185").Value = '%'
$ws.Range("M185").Value = 'T4'

# Step 5: Numeric Value column J
$ws.Range("J177").Value = 245000000
$ws.Range("J178").Value = 205000000
$ws.Range("J179").Value = 360000000
$ws.Range("J180").Value = 869000000
$ws.Range("J181").Value = 863000000
$ws.Range("J182").Value = 850000000
$ws.Range("J183").Value = 17
$ws.Range("J184").Value = 20
$ws.Range("J185").Value = 27

# Step 6: Update the active cell selection to match the saved view state
$ws.Range("H191").Select()
